$d = $word.ActiveDocument

function SplitAt($pos) {
    # Inserting (then deleting) a throw-away bookmark at a zero-length
    # Range is the only edit this engine performs that splits a run
    # in two without re-merging every equally-formatted run nearby.
    $pt = $d.Range($pos, $pos)
    $pt.Bookmarks.Add("TmpSplitMarker")
    $d.Bookmarks("TmpSplitMarker").Delete()
}

# ---------------------------------------------------------------------
# 1) "This jython autopsy module ..." -> "This autopsy module ..."
# ---------------------------------------------------------------------
$r = $d.Content
$r.Start = 0
$r.End = $d.Content.End
$found = $r.Find.Execute("This jython autopsy module can calculate perceptual hash value of jpg files", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $introStart = $r.Start

    # Locate the remaining boundaries (jpg / the long sentence / the comma /
    # " and look for similar pictures.") before anything is edited, they are
    # all relative to the still-untouched original text.
    $jpgR = $d.Content
    $jpgR.Start = 0
    $jpgR.End = $d.Content.End
    [void]$jpgR.Find.Execute("jpg", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $jpgStart = $jpgR.Start
    $jpgEnd = $jpgR.End

    $commaR = $d.Content
    $commaR.Start = 0
    $commaR.End = $d.Content.End
    [void]$commaR.Find.Execute("other pictures' value,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $sentenceEnd = $commaR.End - 1   # just before the comma
    $commaEnd = $commaR.End          # just after the comma (bookmark sits here today)

    # Rewrite the run-1 text in place (this necessarily merges it with its
    # neighbours, since they all share identical/default formatting - the
    # splits below restore the original paragraph-internal run breaks).
    $run1 = $d.Range($introStart, $jpgStart)
    $run1.Text = "This autopsy module can calculate perceptual hash value of "

    $shift = $jpgStart - ($introStart + $run1.Text.Length)

    SplitAt($jpgStart - $shift)
    SplitAt($jpgEnd - $shift)
    SplitAt($sentenceEnd - $shift)
    SplitAt($commaEnd - $shift)
}

# ---------------------------------------------------------------------
# 2) Drop the _GoBack bookmark from its old spot (right after the comma,
#    before " and look for similar pictures.").
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 3) "Line 86: " -> "Line " / "75" / _GoBack bookmark / ": "
# ---------------------------------------------------------------------
$r2 = $d.Content
$r2.Start = 0
$r2.End = $d.Content.End
$found2 = $r2.Find.Execute("Line 86: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found2) {
    $lineStart = $r2.Start
    $lineEnd = $r2.End

    # Change the digits while the paragraph fragment is still one run -
    # editing text, not just splitting, is what triggers the re-merge, so
    # do it before re-establishing any run boundaries.
    $numRange = $d.Range($lineStart + 5, $lineStart + 7)
    $numRange.Text = "75"

    # Recreate "Line " | "75" boundary.
    SplitAt($lineStart + 5)

    # Recreate ": " | "sys.path..." boundary (length preserved, 86->75).
    SplitAt($lineEnd)

    # Insert the real bookmark between "75" and ": " - this also performs
    # that split.
    $bm = $d.Range($lineStart + 7, $lineStart + 7)
    $bm.Bookmarks.Add("_GoBack")
}
